# TC04_Canine_Filter_PrimDisSite-LymphNode.xlsx -- "startup" sheet edit
# Commit: "updated keyword for case files tab data writing issue + more curated scripts"
#
# - CasesTab row: Cases query gains a trailing `Cohort` projection.
# - SamplesTab row: Samples query is untouched (kept for shared-string order).
# - FilesTab row: Files query drops the trailing `Study Code` projection.
# - StatQuery column (C2:C4) is replaced everywhere by a new
#   Programs/Studies/Cases/Samples/Case Files/Study Files summary query.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
WHERE diag.primary_disease_site IN ['Lymph Node']
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
coalesce(co.cohort_description, '') AS `Cohort`
'@

$statQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE diag.primary_disease_site IN ['Lymph Node']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$filesQuery = @'

MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE diag.primary_disease_site IN ['Lymph Node']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis
'@

# Order matters here: the workbook's shared-string table is rebuilt with
# unused strings dropped and new strings appended in assignment order, so
# writing Cases, then StatQuery, then Files (matching the source diff)
# reproduces the same shared-string layout as the authored change.
$ws.Range("B2").Value2 = $casesQuery
$ws.Range("C2").Value2 = $statQuery
$ws.Range("C3").Value2 = $statQuery
$ws.Range("C4").Value2 = $statQuery
$ws.Range("B4").Value2 = $filesQuery

# Rows re-fit to the new wrapped query text (15pt/line at the refreshed
# default row height).
$ws.Rows.Item(2).RowHeight = 270
$ws.Rows.Item(3).RowHeight = 225
$ws.Rows.Item(4).RowHeight = 210

# View state left behind by the editing session.
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("C13:C14").Select()
